$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that already finished with highlight-only (yellow fill, no text) cells: week 8 and week 13
$highlightOnlyRows = @(13, 18)
foreach ($r in $highlightOnlyRows) {
    foreach ($col in @("C", "D", "E")) {
        $ws.Range("$col$r").Interior.Color = 65535
    }
}

# Rows marked "OK" (completed classes): weeks 12, 14, 16, 17, 18, 19, 20
$okRows = @(17, 19, 21, 22, 23, 24, 25)
foreach ($r in $okRows) {
    foreach ($col in @("C", "D", "E")) {
        $ws.Range("$col$r").Value = "OK"
    }
}

# Update the active selection to reflect where the author ended up working
$ws.Range("C27").Select() | Out-Null
